$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: M1, Ccl24, Ccr3, M1 ...
$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Ccl24"
$ws.Range("C2").Value = "Ccr3"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4675276666666666
$ws.Range("H2").Value = 1.402583
$ws.Range("I2").Value = 0.1632394805509195
$ws.Range("J2").Value = 0.1632394805509196
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1323866666666667
$ws.Range("N2").Value = 0.39716
$ws.Range("O2").Value = 0.2977240407890032
$ws.Range("P2").Value = 0.2977240407890032
$ws.Range("Q2").Value = 0.06189442936444443
$ws.Range("R2").Value = 0.5570498642799999
$ws.Range("S2").Value = 0.04860031776591765
$ws.Range("T2").Value = 0.04860031776591766

# Row 3: M1, Ccl24, Ccr3, M2 ...
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Ccl24"
$ws.Range("C3").Value = "Ccr3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4675276666666666
$ws.Range("H3").Value = 1.402583
$ws.Range("I3").Value = 0.1632394805509195
$ws.Range("J3").Value = 0.1632394805509196
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3122756666666667
$ws.Range("N3").Value = 0.936827
$ws.Range("O3").Value = 0.7022759592109968
$ws.Range("P3").Value = 0.7022759592109967
$ws.Range("Q3").Value = 0.1459975137934444
$ws.Range("R3").Value = 1.313977624141
$ws.Range("S3").Value = 0.1146391627850019
$ws.Range("T3").Value = 0.1146391627850019

# Row 4: M2, Ccl24, Ccr3, M1 ...
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Ccl24"
$ws.Range("C4").Value = "Ccr3"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.396532333333333
$ws.Range("H4").Value = 7.189597
$ws.Range("I4").Value = 0.8367605194490804
$ws.Range("J4").Value = 0.8367605194490806
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1323866666666667
$ws.Range("N4").Value = 0.39716
$ws.Range("O4").Value = 0.2977240407890032
$ws.Range("P4").Value = 0.2977240407890032
$ws.Range("Q4").Value = 0.3172689271688889
$ws.Range("R4").Value = 2.85542034452
$ws.Range("S4").Value = 0.2491237230230855
$ws.Range("T4").Value = 0.2491237230230856

# Row 5: M2, Ccl24, Ccr3, M2 ...
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Ccl24"
$ws.Range("C5").Value = "Ccr3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.396532333333333
$ws.Range("H5").Value = 7.189597
$ws.Range("I5").Value = 0.8367605194490804
$ws.Range("J5").Value = 0.8367605194490806
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3122756666666667
$ws.Range("N5").Value = 0.936827
$ws.Range("O5").Value = 0.7022759592109968
$ws.Range("P5").Value = 0.7022759592109967
$ws.Range("Q5").Value = 0.7483787320798889
$ws.Range("R5").Value = 6.735408588718999
$ws.Range("S5").Value = 0.5876367964259949
$ws.Range("T5").Value = 0.587636796425995
